$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (fedfundsrate) - only row 2 changes
$ws.Range("B2").Value = 9.6866666666666656

# Update column C (fedfundsrate_shadow) values
$ws.Range("C2").Value = 9.6866666666666656
$ws.Range("C3").Value = 10.556666666666658
$ws.Range("C4").Value = 11.38999999999999
$ws.Range("C5").Value = 9.2666666666667119
$ws.Range("C7").Value = 7.9233333333333933
$ws.Range("C8").Value = 7.9000000000001069
$ws.Range("C9").Value = 8.1033333333334401
$ws.Range("C10").Value = 7.8266666666667373
$ws.Range("C11").Value = 6.9200000000000372
$ws.Range("C12").Value = 6.2066666666667825
$ws.Range("C13").Value = 6.2666666666667759
$ws.Range("C14").Value = 6.2200000000000921
$ws.Range("C15").Value = 6.650000000000067
$ws.Range("C16").Value = 6.8433333333334678
$ws.Range("C17").Value = 6.9166666666667265
$ws.Range("C18").Value = 6.6633333333334432
$ws.Range("C19").Value = 7.1566666666667667
$ws.Range("C20").Value = 7.9833333333334533
$ws.Range("C21").Value = 8.4700000000001108
$ws.Range("C22").Value = 9.4433333333334701
$ws.Range("C23").Value = 9.7266666666667501
$ws.Range("C24").Value = 9.0833333333333766
$ws.Range("C25").Value = 8.6133333333334292
$ws.Range("C26").Value = 8.250000000000135
$ws.Range("C27").Value = 8.243333333333446
$ws.Range("C28").Value = 8.1600000000001227
$ws.Range("C29").Value = 7.7433333333334353
$ws.Range("C30").Value = 6.4266666666668026
$ws.Range("C31").Value = 5.8633333333333981
$ws.Range("C32").Value = 5.6433333333334001
$ws.Range("C33").Value = 4.8166666666667579
$ws.Range("C34").Value = 4.0233333333334453
$ws.Range("C35").Value = 3.7700000000001177
$ws.Range("C36").Value = 3.2566666666667965
$ws.Range("C37").Value = 3.0366666666667541
$ws.Range("C38").Value = 3.0400000000000871
$ws.Range("C43").Value = 3.9400000000001212
$ws.Range("C47").Value = 6.0200000000000919
$ws.Range("C49").Value = 5.7200000000001028
$ws.Range("C56").Value = 5.5333333333334123
$ws.Range("C81").Value = 0.9963915207263474
$ws.Range("C82").Value = 1.0029171386297664
$ws.Range("C83").Value = 1.0095223043066204
$ws.Range("C84").Value = 1.4328440990858837
$ws.Range("C85").Value = 1.9495334161461209
$ws.Range("C86").Value = 2.4695789344845842
$ws.Range("C87").Value = 2.9429721106922324
$ws.Range("C88").Value = 3.459706093287318
$ws.Range("C89").Value = 3.9797760238860835
$ws.Range("C90").Value = 4.4565115712214887
$ws.Range("C91").Value = 4.9065766149136492
$ws.Range("C92").Value = 5.2466358969100124
$ws.Range("C93").Value = 5.2466879902888897
$ws.Range("C94").Value = 5.2567319878412899
$ws.Range("C95").Value = 5.2501009120170972
$ws.Range("C96").Value = 5.0734612778135624
$ws.Range("C97").Value = 4.4968129986319605
$ws.Range("C98").Value = 3.1768226690157331
$ws.Range("C99").Value = 2.0868252758880823
$ws.Range("C100").Value = 1.9401557219421051
$ws.Range("C101").Value = 0.50775193785634887
$ws.Range("C102").Value = 1.6975116104384069
$ws.Range("C103").Value = 0.67964963396671685
$ws.Range("C104").Value = -0.10241757465069279
$ws.Range("C105").Value = -0.1745308066623652
$ws.Range("C106").Value = -0.16850227418053754
$ws.Range("C107").Value = -1.2577057835605143
$ws.Range("C108").Value = -1.1337212152887521
$ws.Range("C109").Value = -1.5800407870430688
$ws.Range("C110").Value = -1.3608393103700323
$ws.Range("C111").Value = -1.0737038927331355
$ws.Range("C112").Value = -1.8431744468209721
$ws.Range("C113").Value = -1.223111938587762
$ws.Range("C114").Value = -1.9888784185039032
$ws.Range("C115").Value = -1.702175400407191
$ws.Range("C116").Value = -1.3467568371485217
$ws.Range("C117").Value = -2.447145310551746
$ws.Range("C118").Value = -1.1834911429906003
$ws.Range("C119").Value = -0.69560187802735207
$ws.Range("C120").Value = -0.47310380001696073
$ws.Range("C121").Value = -0.56019914395267767
$ws.Range("C122").Value = -0.77095224163514775
$ws.Range("C123").Value = -0.71018212769901057
$ws.Range("C124").Value = -0.45646496664363267
$ws.Range("C125").Value = -0.21272213887576497
$ws.Range("C126").Value = 0.17616224082710463
$ws.Range("C127").Value = 0.036083754907800802
$ws.Range("C128").Value = 0.053799398148735023
$ws.Range("C147").Value = 7.8096896392297888
$ws.Range("C148").Value = -4.7324215093060396
$ws.Range("C149").Value = -3.1385738473021951
$ws.Range("C150").Value = -3.0358164520225661

# Add new row 151
$ws.Range("A151").Value = 2021.25
$ws.Range("B151").Value = 0
$ws.Range("C151").Value = -2.7126338022849694
